# Apply crypto price/symbol list update (GitHub Actions run, 2022-12-29)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Numeric-looking "Price" column (D) values must stay text, like the rest of the sheet ---
# (the sheet stores every cell as inline/shared text -- force text format first so
# Excel does not silently coerce these into real numbers)
$priceUpdates = @{
    "D2" = "246.04"
    "D3" = "24.20"
    "D4" = "5.360"
    "D6" = "6.477"
    "D7" = "3.143"
    "D8" = "0.8189"
    "D9" = "0.8760"
    "D10" = "0.01013"
    "D11" = "0.1380"
    "D12" = "0.06981"
    "D13" = "0.03148"
    "D14" = "0.02942"
    "D15" = "0.09409"
    "D16" = "3.748"
    "D17" = "0.001524"
    "D18" = "0.04710"
    "D19" = "0.006250"
    "D20" = "0.001235"
    "D21" = "0.004687"
    "D22" = "0.00008797"
    "D24" = "2.146"
    "D25" = "0.3181"
    "D40" = "0.03717"
    "D41" = "0.006378"
    "D42" = "0.1060"
    "D43" = "0.002799"
    "D44" = "0.007527"
    "D45" = "0.00005268"
    "D48" = "0.002538"
    "D49" = "0.00002099"
    "D50" = "0.0001999"
}
foreach ($addr in $priceUpdates.Keys) {
    $cell = $ws.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = $priceUpdates[$addr]
}
# Drop back to the default (unstyled) look now that the text is locked in
$ws.Range("D2:D50").Style = "Normal"

# --- Coin / Link / Volume(1h) text columns (plain text, no coercion risk) ---
$textUpdates = @{
    "B10" = "One"
    "C10" = "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
    "E10" = "9OneONEBestin24h"
    "B11" = "WazirX"
    "C11" = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
    "E11" = "10WazirXWRX"
    "B12" = "MandalaExchangeToken"
    "C12" = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
    "E12" = "11MandalaExchangeTokenMDX"
    "B13" = "LiechtensteinCryptoassetsExchange"
    "C13" = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
    "E13" = "12LiechtensteinCryptoassetsExchangeLCX"
    "B14" = "BitrueCoin"
    "C14" = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
    "E14" = "13BitrueCoinBTR"
    "B15" = "BitMartToken"
    "C15" = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
    "E15" = "14BitMartTokenBMX"
    "B16" = "MCDex"
    "C16" = "https://coinranking.com/coin/3nMM61qeg+mcdex-mcb"
    "E16" = "15MCDexMCB"
    "B17" = "BitForexToken"
    "C17" = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
    "E17" = "16BitForexTokenBF"
    "B18" = "CoinExToken"
    "C18" = "https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet"
    "E18" = "17CoinExTokenCET"
    "E43" = "42CEJICEJI"
    "E47" = "46CoinbaseStockTokenCOINWorstin24h"
}
foreach ($addr in $textUpdates.Keys) {
    $ws.Range($addr).Value = $textUpdates[$addr]
}
